$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds plain text dates (e.g. "2025-12-02"), not real Excel dates.
# Temporarily force text format so Excel doesn't auto-convert the string into
# a date serial number, then restore the default "Normal" style so the new
# rows don't pick up a stray explicit style like the original data row.
$dateCells = @("A3", "A4", "A5", "A6", "A7")
foreach ($cell in $dateCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("A3").Value = "2025-12-02"
$ws.Range("B3").Value = "food"
$ws.Range("C3").Value = 100

$ws.Range("A4").Value = "2025-12-02"
$ws.Range("B4").Value = "food"
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = "2025-12-02"
$ws.Range("B5").Value = "soft drinks"
$ws.Range("C5").Value = 100

$ws.Range("A6").Value = "2025-12-05"
$ws.Range("B6").Value = "soft drinks"
$ws.Range("C6").Value = 1

$ws.Range("A7").Value = "2026-02-26"
$ws.Range("B7").Value = "soaps"
$ws.Range("C7").Value = 1

foreach ($cell in $dateCells) {
    $ws.Range($cell).Style = "Normal"
}
